$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that looks like a plain number but must be
# stored as literal text (matching the source data which is always
# written as strings). Using a leading apostrophe forces text entry,
# then the cell style is reset to Normal so no stray formatting
# (quote-prefix) lingers on the cell itself.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.406.84"
$ws.Range("E2").Value = "  -1.91%  "

$ws.Range("D3").Value = "1.792.26"
$ws.Range("E3").Value = "  -2.23%  "

Set-TextValue "D4" "1.008"
$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("E5").Value = "  +0.25%  "

Set-TextValue "D6" "307.52"
$ws.Range("E6").Value = "  -0.91%  "

Set-TextValue "D7" "0.4566"
$ws.Range("E7").Value = "  -1.06%  "

Set-TextValue "D8" "0.3636"
$ws.Range("E8").Value = "  -0.46%  "

Set-TextValue "D9" "46.88"
$ws.Range("E9").Value = "  +1.78%  "

Set-TextValue "D10" "0.07092"
$ws.Range("E10").Value = "  -1.09%  "

Set-TextValue "D11" "0.8758"
$ws.Range("E11").Value = "  -0.32%  "

Set-TextValue "D12" "0.07847"
$ws.Range("E12").Value = "  +0.37%  "

Set-TextValue "D13" "19.55"
$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("D14").Value = "1.817.65"
$ws.Range("E14").Value = "  -1.79%  "

Set-TextValue "D15" "5.278"
$ws.Range("E15").Value = "  -0.88%  "

Set-TextValue "D16" "6.322"
$ws.Range("E16").Value = "  -0.55%  "

Set-TextValue "D17" "85.04"
$ws.Range("E17").Value = "  -3.95%  "

Set-TextValue "D18" "1.009"
$ws.Range("E18").Value = "  +0.28%  "

Set-TextValue "D19" "0.000008534"
$ws.Range("E19").Value = "  -2.44%  "

Set-TextValue "D20" "1.007"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").Value = "26.429.90"
$ws.Range("E21").Value = "  -1.88%  "

Set-TextValue "D22" "14.25"
$ws.Range("E22").Value = "  -1.71%  "

Set-TextValue "D23" "4.989"
$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("D24").Value = "2.074.15"
$ws.Range("E24").Value = "  -0.66%  "

Set-TextValue "D25" "10.51"
$ws.Range("E25").Value = "  +0.89%  "

Set-TextValue "D26" "1.990"
$ws.Range("E26").Value = "  +1.13%  "

Set-TextValue "D27" "152.21"
$ws.Range("E27").Value = "  +1.01%  "

Set-TextValue "D28" "17.91"
$ws.Range("E28").Value = "  -1.60%  "

Set-TextValue "D29" "2.041"
$ws.Range("E29").Value = "  +2.25%  "

Set-TextValue "D30" "112.33"
$ws.Range("E30").Value = "  -1.04%  "

Set-TextValue "D31" "4.859"
$ws.Range("E31").Value = "  -1.71%  "

Set-TextValue "D32" "0.08672"
$ws.Range("E32").Value = "  -1.75%  "

Set-TextValue "D33" "3.054"
$ws.Range("E33").Value = "  -1.56%  "

Set-TextValue "D34" "4.446"
$ws.Range("E34").Value = "  -0.37%  "

Set-TextValue "D35" "0.7255"
$ws.Range("E35").Value = "  -4.33%  "

Set-TextValue "D36" "2.653"
$ws.Range("E36").Value = "  +0.11%  "

Set-TextValue "D37" "1.108"
$ws.Range("E37").Value = "  -2.47%  "

Set-TextValue "D38" "1.006"
$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("E39").Value = "  -1.32%  "

Set-TextValue "D40" "0.01944"
$ws.Range("E40").Value = "  +1.03%  "

$ws.Range("E41").Value = "  -0.29%  "

Set-TextValue "D42" "0.5276"
$ws.Range("E42").Value = "  +6.00%  "

Set-TextValue "D43" "2.870"
$ws.Range("E43").Value = "  -1.75%  "

Set-TextValue "D44" "6.905"
$ws.Range("E44").Value = "  -0.48%  "

Set-TextValue "D45" "0.1516"
$ws.Range("E45").Value = "  -4.98%  "

Set-TextValue "D46" "8.022"
$ws.Range("E46").Value = "  -4.04%  "

Set-TextValue "D47" "0.4722"
$ws.Range("E47").Value = "  +1.11%  "

$ws.Range("E48").Value = "  +0.28%  "

Set-TextValue "D49" "9.855"
$ws.Range("E49").Value = "  -3.25%  "

Set-TextValue "D50" "99.92"
$ws.Range("E50").Value = "  -2.60%  "

$ws.Range("E51").Value = "  -1.27%  "
